$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B9").Value = 4477
$ws.Range("C9").Value = 4238
$ws.Range("D9").Value = 4454
$ws.Range("E9").Value = 4476
$ws.Range("F9").Value = 4477
